# iteration 4 - remise final
#
# Applies the "iteration 4" update to estimation.xlsx:
#  - fills in the Iteration #3 (sheet index 4) wrap-up fields (B40/B42)
#  - fills in the Iteration #4 (sheet index 5) task log (rows 14-17),
#    the total (C37, recalculated), and the wrap-up fields (B40/B42)
#  - moves the active tab / selection from Iteration #3 to Iteration #4

$wb = $excel.ActiveWorkbook

$wsIter2 = $wb.Worksheets.Item(3)   # "Iteration #2"
$wsIter3 = $wb.Worksheets.Item(4)   # "Iteration #3"
$wsIter4 = $wb.Worksheets.Item(5)   # "Iteration #4"

# ---------------------------------------------------------------------
# Iteration #3 sheet: a stray date-style cell + the self-eval wrap-up
# ---------------------------------------------------------------------

# A21 picks up the same "d-mmm" date style used by the other date cells
# in the column (no value is entered, just the formatting).
$wsIter3.Range("A21").NumberFormat = "d-mmm"

# Self-evaluation note + grade for iteration #3
$wsIter3.Range("B40").Value = 6
$wsIter3.Range("B42").Value = "beaucoup de correction appliqués dans cette itération, ajout dune base de donnée"

# ---------------------------------------------------------------------
# Iteration #4 sheet: new task-log rows + wrap-up
# ---------------------------------------------------------------------

# Row 14
$wsIter4.Range("A14").Value = "2018-05-06"
$wsIter4.Range("B14").Value = "ajout taux protection, détails dans ressources. Création de l'activity ""Survivant"", tentative de création d'un survivant et lajouter dans une liste ensuite et afficher cette liste dans une listView"
# entered as text (quote-prefixed) rather than the number 7.5
$wsIter4.Range("C14").Formula = "'7.5"

# Row 15
$wsIter4.Range("A15").NumberFormat = "d-mmm"
$wsIter4.Range("A15").Value = "2018-05-07"
$wsIter4.Range("B15").Value = "gestion de lactivité survivant"
# entered as text (quote-prefixed) rather than the number 4.5
$wsIter4.Range("C15").Formula = "'4.5"

# Row 16
$wsIter4.Range("A16").NumberFormat = "d-mmm"
$wsIter4.Range("A16").Value = "2018-05-08"
$wsIter4.Range("B16").Value = "gestion de la suppression dun survivant "
$wsIter4.Range("C16").Value = 3

# Row 17
$wsIter4.Range("A17").NumberFormat = "d-mmm"
$wsIter4.Range("A17").Value = "2018-05-14"
$wsIter4.Range("B17").Value = "tentative correction supprimerUser dans activité"
$wsIter4.Range("C17").Value = 1

# Self-evaluation note + grade for iteration #4
$wsIter4.Range("B40").Value = 8
$wsIter4.Range("B42").Value = "ajout de lactivity survivant qui permet dafficher les données des survivants. Gestion des relative/linear layout mieux executer car plus de compréhension qu'au début du projet. "

# Recalculate so the C37 = SUM(C14:C36) total picks up the new rows
$excel.Calculate()

# ---------------------------------------------------------------------
# Window / selection state: move focus from Iteration #3 to Iteration #4
# ---------------------------------------------------------------------

$wsIter2.Range("B42:B47").Select()
$wsIter3.Range("B42:B47").Select()
$wsIter4.Range("B42:B47").Select()
$wsIter4.Activate()

# Scroll the sheet tab strip so tab 2 ("Iteration #1") is the first visible
# tab (best effort - matches firstSheet="1" on the workbook view).
$excel.ActiveWindow.ScrollWorkbookTabs(1)
